$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill A2 and A3 with the same text as A1 (shared string "version6")
$ws.Range("A2").Value = "version6"
$ws.Range("A3").Value = "version6"

# Move the active selection to A4, as in the edited workbook
$ws.Range("A4").Select()
